$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the one changed description (row 33, column M / productAriaLabel)
$ws.Cells.Item(33, 13).Value = "Prix Garantie Frischback Baguettes - Online kein Bestand 1.00 Schweizer Franken"

# 2. Insert a new row at 400 (pushes old row 400 down to 401) and fill it in
#    with the new "Motta Panettone" product that was scraped in between.
$ws.Rows.Item(400).Insert()

# Columns whose text looks like a pure number ("3351706", "2.70", "3.00")
# get auto-coerced to a Number by the Value setter, same as typing them into
# Excel would. Force the cell to Text first so the literal string is kept,
# then drop the number-format override so no stray style sticks around.
$ws.Cells.Item(400, 1).NumberFormat = "@"
$ws.Cells.Item(400, 1).Value = "3351706"
$ws.Cells.Item(400, 1).ClearFormats()

$ws.Cells.Item(400, 2).Value = "Motta Panettone"
$ws.Cells.Item(400, 3).Value = "/de/lebensmittel/suesses-snacks/guetzli-suessgebaeck/suesse-broetchen-panettone/motta-panettone/p/3351706"
$ws.Cells.Item(400, 4).Value = 1
$ws.Cells.Item(400, 5).Value = 5
$ws.Cells.Item(400, 6).Value = "Motta"

$ws.Cells.Item(400, 7).NumberFormat = "@"
$ws.Cells.Item(400, 7).Value = "2.70"
$ws.Cells.Item(400, 7).ClearFormats()

$ws.Cells.Item(400, 8).Value = "3.00/100g"
$ws.Cells.Item(400, 9).Value = "Preis pro 100 Gramm"

$ws.Cells.Item(400, 10).NumberFormat = "@"
$ws.Cells.Item(400, 10).Value = "3.00"
$ws.Cells.Item(400, 10).ClearFormats()

$ws.Cells.Item(400, 11).Value = "100g"
$ws.Cells.Item(400, 12).Value = "['lebensmittel', 'suesses-snacks', 'guetzli-suessgebaeck', 'suesse-broetchen-panettone']"
$ws.Cells.Item(400, 13).Value = "Motta Panettone 2.70 Schweizer Franken"

# declarationIcons (N) is blank for every row; materialise an empty cell so
# row 400 keeps the same 15-column shape as all the others.
$ws.Cells.Item(400, 14).NumberFormat = "General"
$ws.Cells.Item(400, 14).ClearFormats()

$ws.Cells.Item(400, 15).Value = "2023-01-03 12:56:48"

# 3. Refresh the scrape timestamp (column O) for every data row, 2..401
#    (401 now, since the insert above shifted everything from 400 down by one).
for ($i = 2; $i -le 401; $i++) {
    $ws.Cells.Item($i, 15).Value = "2023-01-03 12:56:48"
}
